# plotEIC methods for fGroupsSet
#
# 1. Insert a new row above row 19 ("getFeatures") for the new function
#    "getEICsForFGroups" (implement=X, ionize=X, done=X), pushing all
#    subsequent rows down by one.
# 2. Mark the existing "plotEIC" row (now row 34 after the shift) as done
#    in column G.
# 3. Update the active selection to G35 (matches the new sheet state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 19, shifting rows 19:53 down to 20:54.
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with the new function entry.
$ws.Range("A19").Value = "getEICsForFGroups"
$ws.Range("D19").Value = "X"
$ws.Range("F19").Value = "X"
$ws.Range("G19").Value = "X"

# plotEIC is now on row 34 (was row 33) - mark it as done.
$ws.Range("G34").Value = "X"

# Restore the view's active cell/selection.
$ws.Range("G35").Select()
